# Dev IV Project Rubric.xlsx — applies the commit:
# "Added bloom, cascaded shadow maps, fixed performance on window resize"
#
# The rubric sheet lists individual features in column A; columns B/C/D hold
# the point value of that feature for Milestone I/II/III respectively, and a
# completed feature is marked by putting the milestone numeral ("I","II", or
# "III") in column E and an "X" in column F.  This change records that three
# already-listed Milestone II features were completed (window-resize camera
# fix + full-screen post process + black&white post process + post-process
# chain), adds Bloom as a completed Milestone II feature, and turns the old
# "you can add teacher approved features" placeholder row into a real,
# completed "Cascaded Shadow Maps" feature row (also Milestone II). It also
# flags the two "effective project quality" rows (GIT usage / cleaned up API
# objects) complete for Milestone I.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Row 17: blank placeholder row becomes a real, completed feature ---
$ws.Range("A17").Value = "Cascaded Shadow Maps"
$ws.Range("B17").Value = 4
$ws.Range("C17").Value = 4
$ws.Range("D17").Value = 4
$ws.Range("E17").Value = "II"
$ws.Range("F17").Value = "X"

# --- Row 24: Infinite Sky Box -> completed on Milestone II ---
$ws.Range("E24").Value = "II"
$ws.Range("F24").Value = "X"

# --- Row 57: Camera position/aspect ratio preserved on resize -> completed on Milestone II ---
$ws.Range("E57").Value = "II"
$ws.Range("F57").Value = "X"

# --- Row 72: Full screen post process (NDC space quad) -> completed on Milestone II ---
$ws.Range("E72").Value = "II"
$ws.Range("F72").Value = "X"

# --- Row 73: Black and White/Sepia Scene -> completed on Milestone II ---
$ws.Range("E73").Value = "II"
$ws.Range("F73").Value = "X"

# --- Row 75: Bloom Oversaturation(Glow) Effect -> completed on Milestone II ---
$ws.Range("E75").Value = "II"
$ws.Range("F75").Value = "X"

# --- Row 78: Post Processing System (Chain) -> completed on Milestone II ---
$ws.Range("E78").Value = "II"
$ws.Range("F78").Value = "X"

# --- Effective project quality: mark GIT usage + cleaned up API objects complete ---
$ws.Range("C91").Value = "X"
$ws.Range("C92").Value = "X"

# --- Selection moved from F35 to E25 ---
$ws.Range("E25").Select()

$wb.Application.Calculate()
